$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'287.13"
$ws.Range("E2").Value = "'1.15%"
$ws.Range("E3").Value = "'2.63%"
$ws.Range("D4").Value = "'5.207"
$ws.Range("E4").Value = "'2.07%"
$ws.Range("D5").Value = "'0.06997"
$ws.Range("E5").Value = "'5.67%"
$ws.Range("D6").Value = "'7.428"
$ws.Range("E6").Value = "'1.85%"
$ws.Range("D7").Value = "'3.558"
$ws.Range("E7").Value = "'5.47%"
$ws.Range("D8").Value = "'1.409"
$ws.Range("E8").Value = "'3.79%"
$ws.Range("D9").Value = "'0.9017"
$ws.Range("E9").Value = "'-3.71%"
$ws.Range("D10").Value = "'0.1607"
$ws.Range("E10").Value = "'2.39%"
$ws.Range("D11").Value = "'0.07534"
$ws.Range("E11").Value = "'24.47%"
$ws.Range("D12").Value = "'0.07676"
$ws.Range("E12").Value = "'1.66%"
$ws.Range("D13").Value = "'0.02931"
$ws.Range("E13").Value = "'1.64%"
$ws.Range("D14").Value = "'0.09016"
$ws.Range("E14").Value = "'0.89%"
$ws.Range("D15").Value = "'0.001587"
$ws.Range("E15").Value = "'-0.57%"
$ws.Range("D16").Value = "'0.0006528"
$ws.Range("E16").Value = "'1.44%"
$ws.Range("D17").Value = "'0.006186"
$ws.Range("E17").Value = "'-1.51%"
$ws.Range("D18").Value = "'3.473"
$ws.Range("E18").Value = "'0.01%"
$ws.Range("D19").Value = "'2.228"
$ws.Range("E19").Value = "'-0.49%"
$ws.Range("D20").Value = "'0.3239"
$ws.Range("E20").Value = "'1.45%"
$ws.Range("D21").Value = "'0.1333"
$ws.Range("E21").Value = "'2.45%"
$ws.Range("D22").Value = "'4.014"
$ws.Range("E22").Value = "'-1.62%"
$ws.Range("D23").Value = "'0.1599"
$ws.Range("E23").Value = "'5.48%"
$ws.Range("E24").Value = "'1.34%"
$ws.Range("D25").Value = "'0.001209"
$ws.Range("E25").Value = "'2.87%"
$ws.Range("D26").Value = "'0.004242"
$ws.Range("E26").Value = "'-4.83%"
$ws.Range("D27").Value = "'0.0001169"
$ws.Range("E27").Value = "'-6.26%"
$ws.Range("D28").Value = "'0.0001668"
$ws.Range("E28").Value = "'3.70%"
$ws.Range("D40").Value = "'0.04345"
$ws.Range("E40").Value = "'4.37%"
$ws.Range("D41").Value = "'0.006946"
$ws.Range("E41").Value = "'5.20%"
$ws.Range("E42").Value = "'-0.13%"
$ws.Range("D43").Value = "'0.002068"
$ws.Range("E43").Value = "'2.63%"
$ws.Range("D44").Value = "'0.01160"
$ws.Range("E44").Value = "'1.10%"
$ws.Range("D45").Value = "'0.00005843"
$ws.Range("E45").Value = "'6.26%"
$ws.Range("D47").Value = "'0.01306"
$ws.Range("E47").Value = "'0.53%"